$d = $word.ActiveDocument

# --- Edit 1: first paragraph -------------------------------------------
# "This is a Microsoft word document." -> same text + two trailing spaces,
# followed by a red-colored parenthetical note, typed/applied as three
# separate runs (as happened in the real edit).
$d.Content.Find.Execute("This is a Microsoft word document.", $true, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$para1 = $d.Paragraphs(1).Range
$ip1 = $para1.End - 1
$r1 = $d.Range($ip1, $ip1)
$r1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 255

$para1 = $d.Paragraphs(1).Range
$ip2 = $para1.End - 1
$r2 = $d.Range($ip2, $ip2)
$r2.InsertAfter("rsion for main branch")
$r2.Font.Color = 255

$para1 = $d.Paragraphs(1).Range
$ip3 = $para1.End - 1
$r3 = $d.Range($ip3, $ip3)
$r3.InsertAfter(")")
$r3.Font.Color = 255

# --- Edit 2: drop the trailing "ank God almighty, we are free at last."
# paragraph (the very last paragraph in the document body).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Delete()

# --- Edit 3: that paragraph was the only user of a bunch of now-unused
# styles; Word prunes them from styles.xml. Delete highest-original-index
# first so the by-name lookup (resolved against the original ordinal)
# doesn't run past the shrinking collection.
$unusedStyles = @(
  "podcast-tools__subscribe-links",
  "generic-title",
  "subscribe-more-info",
  "subscribe",
  "audio-tool",
  "Heading 4 Char",
  "Heading 2 Char",
  "Hyperlink",
  "apple-converted-space",
  "Heading 4",
  "Heading 2"
)
foreach ($styleName in $unusedStyles) {
  $style = $d.Styles($styleName)
  $style.Delete()
}
